# Angular Beginner.xlsx - "JS first 10 topics added in Technology Topics"
#
# Adds a new "L" column of Technology-Topic labels alongside the existing
# checklist rows on the "Headings" sheet (rows 31-45), tweaks the B4 value
# on the "Spent" sheet, and updates the active-sheet / selection state to
# match where the author left off editing.

$wb = $excel.ActiveWorkbook
$headings = $wb.Worksheets.Item("Headings")
$spent = $wb.Worksheets.Item("Spent")

# --- New "Technology Topics" column (L) on the Headings sheet ---------
$topics = @{
    31 = "LifecycleHook - ( ngDoCheck, ngAfterViewInit)"
    32 = "LifeCycleHook - (ngAfterViewChecked to destory)"
    33 = "Template Reference variable, ngTemplate, ngContainer, ngtemplateOutet"
    34 = "@HostListener and Host binding"
    35 = "@Input, @Output, @Eventemitter"
    36 = "Pass data from parent to child | @Input decorator"
    37 = "Pass data from child to parent | @Output decorator"
    38 = "Pass data between slibings component | Observable"
    39 = "Observable"
    40 = "Rxjs Subject"
    41 = "pipes | Inbuilt pipes | Custom pipes"
    42 = "Impure Pipe | Async Pipe"
    43 = "Routing & Navigation"
    44 = "Service, Dependency Injection"
    45 = "HTTP service | curd operation"
}

foreach ($row in 31..45) {
    $headings.Range("L$row").Value = $topics[$row]
}

# --- Spent sheet value tweak -------------------------------------------
$spent.Range("B4").Value = 0.3

# --- Selection / active sheet state -------------------------------------
# Set the (inactive) Spent sheet's remembered selection first, then finish
# on the Headings sheet so it ends up the active/selected tab - matching
# Excel's behaviour where selecting a range implicitly activates its sheet.
$spent.Range("I15").Select()

$headings.Activate()
$headings.Range("C30").Select()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
